$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The canonical edit inserts 3 new data rows right before the existing
# row 1140, pushing the previous rows 1140-1238 down to 1143-1241 and
# growing the sheet's used range from A1:R1238 to A1:R1241.
$ws.Range("A1140:A1142").EntireRow.Insert()

function Set-DataRow($row, $fecha, $codreg, $catId, $categoria, $variedad, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $unidad, $origen, $precioKg, $kgUnidades, $clasificacion) {
    $ws.Cells.Item($row, 1).Value2 = 5
    $ws.Cells.Item($row, 2).Value2 = "Macroferia Regional de Talca"
    $ws.Cells.Item($row, 3).Value2 = "Maule"
    $ws.Cells.Item($row, 4).Value2 = $fecha
    $ws.Cells.Item($row, 5).Value2 = $codreg
    $ws.Cells.Item($row, 6).Value2 = $catId
    $ws.Cells.Item($row, 7).Value2 = $categoria
    $ws.Cells.Item($row, 8).Value2 = $variedad
    $ws.Cells.Item($row, 9).Value2 = $calidad
    $ws.Cells.Item($row, 10).Value2 = $volumen
    $ws.Cells.Item($row, 11).Value2 = $precioMin
    $ws.Cells.Item($row, 12).Value2 = $precioMax
    $ws.Cells.Item($row, 13).Value2 = $precioProm
    $ws.Cells.Item($row, 14).Value2 = $unidad
    $ws.Cells.Item($row, 15).Value2 = $origen
    $ws.Cells.Item($row, 16).Value2 = $precioKg
    $ws.Cells.Item($row, 17).Value2 = $kgUnidades
    $ws.Cells.Item($row, 18).Value2 = $clasificacion
}

Set-DataRow 1140 45223 7 100112020 "Tomate" "Larga vida" "Primera" 2500 19000 19000 19000 "`$/bandeja 18 kilos" "Región de Arica y Parinacota" 1056 18 "Hortaliza"

Set-DataRow 1141 45223 7 100112020 "Tomate" "Larga vida" "Primera" 1000 23000 23000 23000 "`$/bandeja 18 kilos" "Región del Maule" 1278 18 "Hortaliza"

Set-DataRow 1142 45223 7 100112020 "Tomate" "Larga vida" "Segunda" 1000 15000 15000 15000 "`$/bandeja 18 kilos" "Región de Arica y Parinacota" 833 18 "Hortaliza"
